# fix: alterar python version para 3.11.5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=59174; B="Danilo Fogaça";        C="TI";                     D="Problemas pessoais"; E=1; F=45104; G=9268.99}
    @{Row=3;  A=70753; B="Clarice Pereira";       C="Engenharia";             D="Consulta médica";    E=2; F=45088; G=6224.76}
    @{Row=4;  A=19618; B="Julia Azevedo";         C="P&D";                    D="Viagem de negócios";  E=6; F=45085; G=4325.74}
    @{Row=5;  A=77537; B="Laís da Costa";         C="Marketing";              D="Problemas pessoais"; E=7; F=45092; G=10038.7}
    @{Row=6;  A=17864; B="João Guilherme Costa";  C="Marketing";              D="Outros";              E=7; F=45082; G=8098.48}
    @{Row=7;  A=12544; B="Ana Clara Lopes";       C="P&D";                    D="Outros";              E=6; F=45096; G=10535.36}
    @{Row=8;  A=19664; B="Henrique da Rocha";     C="Recursos Humanos";       D="Consulta médica";    E=7; F=45086; G=10106.72}
    @{Row=9;  A=10936; B="Enzo da Rocha";         C="Atendimento ao Cliente"; D="Outros";              E=7; F=45100; G=11410.21}
    @{Row=10; A=48670; B="Davi Lucca Aragão";     C="Financeiro";             D="Consulta médica";    E=5; F=45099; G=6240.23}
    @{Row=11; A=80387; B="Luigi Caldeira";        C="Jurídico";               D="Problemas pessoais"; E=6; F=45083; G=5444.49}
)

foreach ($rowData in $data) {
    $r = $rowData.Row
    $ws.Cells.Item($r, 1).Value = $rowData.A
    $ws.Cells.Item($r, 2).Value = $rowData.B
    $ws.Cells.Item($r, 3).Value = $rowData.C
    $ws.Cells.Item($r, 4).Value = $rowData.D
    $ws.Cells.Item($r, 5).Value = $rowData.E
    $ws.Cells.Item($r, 6).Value = $rowData.F
    $ws.Cells.Item($r, 7).Value = $rowData.G
}
